$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) for every data row (2-16) moves from 46063 to 46064.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = 46064
}

# Rows 6 and 8-16 are re-sorted. The content (Beteckning/A, Datum/B, Area/G)
# of each row is replaced with the content that ends up there after sorting.
# (All other columns for these rows share identical template values, so only
# A, B and G need to be rewritten.)

# row 6 <- old row 8
$ws.Cells.Item(6, 1).Value = "A 2593-2024"
$ws.Cells.Item(6, 2).Value = 45313.69204861111
$ws.Cells.Item(6, 7).Value = 2.3

# row 8 <- old row 6
$ws.Cells.Item(8, 1).Value = "A 5792-2024"
$ws.Cells.Item(8, 2).Value = 45335
$ws.Cells.Item(8, 7).Value = 5.6

# row 9 <- old row 14
$ws.Cells.Item(9, 1).Value = "A 13651-2023"
$ws.Cells.Item(9, 2).Value = 45006
$ws.Cells.Item(9, 7).Value = 2.2

# row 10 <- old row 12
$ws.Cells.Item(10, 1).Value = "A 8194-2025"
$ws.Cells.Item(10, 2).Value = 45708
$ws.Cells.Item(10, 7).Value = 1.9

# row 11 <- old row 13
$ws.Cells.Item(11, 1).Value = "A 50997-2025"
$ws.Cells.Item(11, 2).Value = 45946
$ws.Cells.Item(11, 7).Value = 1.5

# row 12 <- old row 10
$ws.Cells.Item(12, 1).Value = "A 35642-2023"
$ws.Cells.Item(12, 2).Value = 45147
$ws.Cells.Item(12, 7).Value = 1.2

# row 13 <- old row 16
$ws.Cells.Item(13, 1).Value = "A 7827-2026"
$ws.Cells.Item(13, 2).Value = 46062.63958333333
$ws.Cells.Item(13, 7).Value = 2.1

# row 14 <- old row 15
$ws.Cells.Item(14, 1).Value = "A 7814-2026"
$ws.Cells.Item(14, 2).Value = 46062.61388888889
$ws.Cells.Item(14, 7).Value = 1.1

# row 15 <- old row 11
$ws.Cells.Item(15, 1).Value = "A 28288-2023"
$ws.Cells.Item(15, 2).Value = 45099.6349537037
$ws.Cells.Item(15, 7).Value = 0.5

# row 16 <- old row 9
$ws.Cells.Item(16, 1).Value = "A 7333-2025"
$ws.Cells.Item(16, 2).Value = 45703.35899305555
$ws.Cells.Item(16, 7).Value = 0.9
